$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Datetime" column (D) used to be filled in with the literal text
# "None" for most rows. Remove that placeholder text everywhere in the
# table body (rows 2-7), but keep/apply the underlined "empty" look that
# row 5 already had - i.e. leave the cell present & formatted, just blank.
foreach ($r in 2..7) {
    $cell = $ws.Range("D$r")
    $cell.Value = ""
    $cell.Font.Underline = $true
}

# Column widths so Task/Datetime comfortably fit their content.
$ws.Columns.Item(3).ColumnWidth = 15
$ws.Columns.Item(4).ColumnWidth = 8.5

# Reflect the new selection over the whole Datetime column's data.
$ws.Range("D2:D7").Select() | Out-Null
